# Weekly update: insert a new price observation row for
# "Feria Lagunitas de Puerto Montt - Acelga" dated 44827, pushing the
# existing rows 124..221 down to 125..222 (dimension grows to A1:R222).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 124; all rows below shift down by one.
$ws.Rows.Item(124).Insert()

# Populate the newly inserted row 124 with the new observation.
$ws.Cells.Item(124, 1).Value = 4
$ws.Cells.Item(124, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(124, 3).Value = "Los Lagos"
$ws.Cells.Item(124, 4).Value = 44827
$ws.Cells.Item(124, 5).Value = 10
$ws.Cells.Item(124, 6).Value = 100112009
$ws.Cells.Item(124, 7).Value = "Acelga"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 200
$ws.Cells.Item(124, 11).Value = 1500
$ws.Cells.Item(124, 12).Value = 1500
$ws.Cells.Item(124, 13).Value = 1500
$ws.Cells.Item(124, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(124, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(124, 16).Value = 1000
$ws.Cells.Item(124, 17).Value = 1.5
$ws.Cells.Item(124, 18).Value = "Hortaliza"
